$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A58").Value = 57
$ws.Range("B58").Value = 1
$ws.Range("C58").Value = "2024-06-16 03:14:08"
$ws.Range("D58").Value = 200
$ws.Range("E58").Value = 3

$ws.Range("A59").Value = 58
$ws.Range("B59").Value = 2
$ws.Range("C59").Value = "2024-06-16 03:14:08"
$ws.Range("D59").Value = 200
$ws.Range("E59").Value = 0
